$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data is being re-assigned (rotated) across rows 2-5 and 7-13; row 6 stays unchanged.

# Row 2
$ws.Range("D2").Value = 44662
$ws.Range("K2").Value = 'Sin especificar'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44662
$ws.Range("K3").Value = 'Sin especificar'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 889
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44280
$ws.Range("K4").Value = 'Sin especificar'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 360000
$ws.Range("O4").Value = 360000
$ws.Range("P4").Value = 360000
$ws.Range("Q4").Value = '$/bins (450 kilos)'
$ws.Range("R4").Value = 'Provincia del Elquí'
$ws.Range("S4").Value = 800
$ws.Range("T4").Value = 450

# Row 5
$ws.Range("D5").Value = 44285
$ws.Range("K5").Value = 'Wonderfull'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 8
$ws.Range("N5").Value = 280000
$ws.Range("O5").Value = 300000
$ws.Range("P5").Value = 290000
$ws.Range("Q5").Value = '$/bins (400 kilos)'
$ws.Range("R5").Value = 'Provincia del Elquí'
$ws.Range("S5").Value = 725
$ws.Range("T5").Value = 400

# Row 7
$ws.Range("D7").Value = 44312
$ws.Range("K7").Value = 'Wonderfull'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 24
$ws.Range("N7").Value = 220000
$ws.Range("O7").Value = 240000
$ws.Range("P7").Value = 230000
$ws.Range("Q7").Value = '$/bins (400 kilos)'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 575
$ws.Range("T7").Value = 400

# Row 8
$ws.Range("D8").Value = 44312
$ws.Range("K8").Value = 'Wonderfull'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 34
$ws.Range("N8").Value = 240000
$ws.Range("O8").Value = 240000
$ws.Range("P8").Value = 240000
$ws.Range("Q8").Value = '$/bins (450 kilos)'
$ws.Range("R8").Value = 'Provincia del Elquí'
$ws.Range("S8").Value = 533
$ws.Range("T8").Value = 450

# Row 9
$ws.Range("D9").Value = 44320
$ws.Range("K9").Value = 'Wonderfull'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 12
$ws.Range("N9").Value = 250000
$ws.Range("O9").Value = 260000
$ws.Range("P9").Value = 255000
$ws.Range("Q9").Value = '$/bins (400 kilos)'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 638
$ws.Range("T9").Value = 400

# Row 10
$ws.Range("D10").Value = 44334
$ws.Range("K10").Value = 'Wonderfull'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 16
$ws.Range("N10").Value = 240000
$ws.Range("O10").Value = 250000
$ws.Range("P10").Value = 245000
$ws.Range("Q10").Value = '$/bins (450 kilos)'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 544
$ws.Range("T10").Value = 450

# Row 11
$ws.Range("D11").Value = 44721
$ws.Range("K11").Value = 'Wonderfull'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 7
$ws.Range("N11").Value = 300000
$ws.Range("O11").Value = 300000
$ws.Range("P11").Value = 300000
$ws.Range("Q11").Value = '$/bins (400 kilos)'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 750
$ws.Range("T11").Value = 400

# Row 12
$ws.Range("D12").Value = 44266
$ws.Range("K12").Value = 'Wonderfull'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 4800
$ws.Range("O12").Value = 4800
$ws.Range("P12").Value = 4800
$ws.Range("Q12").Value = '$/bandeja 4 kilos'
$ws.Range("R12").Value = 'Provincia del Elquí'
$ws.Range("S12").Value = 1200
$ws.Range("T12").Value = 4

# Row 13
$ws.Range("D13").Value = 44266
$ws.Range("K13").Value = 'Wonderfull'
$ws.Range("L13").Value = 'Tercera'
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 4000
$ws.Range("O13").Value = 4000
$ws.Range("P13").Value = 4000
$ws.Range("Q13").Value = '$/bandeja 4 kilos'
$ws.Range("R13").Value = 'Provincia del Elquí'
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 4
